# Journal de travail AGI - update: add "Ajout de la page commande à la
# maquette et ajout de use case scénarios" entry in column A of row 8
# (same wrapped-text style as the existing long entries in column A),
# and grow row 8 to fit the two-line wrapped text (matches row 7's 30pt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 1).Value = "Ajout de la page commande à la maquette et ajout de use case scénarios"
$ws.Cells.Item(8, 1).WrapText = $true
$ws.Rows.Item(8).RowHeight = 30
